$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B4").Copy()
$ws.Range("B5:B7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B5").Value = 45441
$ws.Range("C5").Value = 3
$ws.Range("D5").Value = "Backend object classes, REST endpoints"

$ws.Range("B6").Value = 45442
$ws.Range("C6").Value = 5
$ws.Range("D6").Value = "Backend object model assember, ThunderClient testing"

$ws.Range("B7").Value = 45443
$ws.Range("C7").Value = 2
$ws.Range("D7").Value = "Backend test dataLoader thinkering"

$ws.Range("D8").Select()
